$wb = $excel.ActiveWorkbook

# --- Sheet "Requerimientos": mark R14 ("Subir de nivel") as done ---
$wsReq = $wb.Worksheets.Item("Requerimientos")
$wsReq.Range("B14").Value = "x"

# --- Sheet "Must Have": add new sub-task note for the "Pelotas" row ---
$wsMust = $wb.Worksheets.Item("Must Have")
$wsMust.Activate()
# Copy formatting from the neighbouring filled-in cell (C5) onto F5, then set its text
$wsMust.Range("C5").Copy() | Out-Null
$wsMust.Range("F5").PasteSpecial(-4122) | Out-Null
$wsMust.Range("F5").Value = "Ver si hay pelotas vivas"
$excel.ActiveWindow.ScrollColumn = 2
$wsMust.Range("F5").Select() | Out-Null

# Leave "Requerimientos" as the active/selected tab, with E14 selected
$wsReq.Activate()
$wsReq.Range("E14").Select() | Out-Null
